$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder rows 3-8 (A:B) so the Bigs/Littles pairings match the new layout.
# Row 1,2,4,9 are unchanged.
$ws.Range("A3").Value = "David Zhao*"
$ws.Range("B3").Value = "Ernie and Bert Sanders"

$ws.Range("A5").Value = "Disha Jain"

$ws.Range("A6").Value = "Robyn Guarriello"
$ws.Range("B6").Value = "Ben Cars- My Luggage"

$ws.Range("A7").Value = "Rohni Awasthi"
$ws.Range("B7").Value = "Larry Richards"

$ws.Range("B8").Value = "Ron Ferretly"
